$d = $word.ActiveDocument

# Locate the target paragraph ("The execution times of both insert and
# append...") robustly via Find rather than a hard-coded paragraph index.
$findRange = $d.Content
$found = $findRange.Find.Execute(
    "The insert function grows much faster with the input size than the append function.",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if (-not $found) {
    throw "Could not locate the target paragraph."
}

$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($findRange.Start -ge $p.Range.Start -and $findRange.Start -lt $p.Range.End) {
        $targetPara = $p
    }
}
if ($null -eq $targetPara) {
    throw "Could not resolve the paragraph object for the found text."
}

# Collapse a copy of its range to the very end (just before the paragraph
# mark) and inject raw OOXML: the original paragraph re-asserted (so its
# w14:paraId/rsid attributes survive), a new blank paragraph, and the new
# paragraph explaining Array.unshift (with spell-check proofErr markers
# around the non-dictionary "Array.unshift" token, matching what Word's
# proofer would emit for freshly typed text).
$insertionPoint = $targetPara.Range.Duplicate
$insertionPoint.Collapse(0)

$xml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml" w14:paraId="25CA273C" w14:textId="574FF333" w:rsidR="00185803" w:rsidRDefault="000D0FBD"><w:r><w:t>The execution times of both insert and append grow as the input size grows. The insert function grows much faster with the input size than the append function.</w:t></w:r></w:p>' + `
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>' + `
    '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:r><w:t xml:space="preserve">The insert function uses the </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>Array.unshift</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> method, which has a time complexity of O(n). It has to increment each element in the existing array by 1 to add the new element to the beginning of the array.</w:t></w:r></w:p>'

$insertionPoint.InsertXML($xml)

Write-Output "Inserted Array.unshift explanation paragraph. Paragraph count now $($d.Paragraphs.Count)."
